# Fixar (pin) o banco de dados de autoria e mencoes no ano de analise.
#
# As colunas de "qtd" (sheet 1) e "tx-sucesso" (sheet 5) tem UFs empatadas
# (mesmo valor numerico) cuja ordem de desempate estava instavel. Este
# script fixa a ordem de exibicao das UFs empatadas, sem alterar nenhum
# valor numerico.

$wb = $excel.ActiveWorkbook

# --- aba "qtd": troca PB/BA (empate em 30) e rotaciona RO/AM/MA (empate em 4)
$wsQtd = $wb.Worksheets.Item("qtd")
$wsQtd.Range("A10").Value = "BA"
$wsQtd.Range("A11").Value = "PB"
$wsQtd.Range("A20").Value = "MA"
$wsQtd.Range("A21").Value = "RO"
$wsQtd.Range("A22").Value = "AM"

# --- aba "tx-sucesso": rotaciona MT/MA/XX (empate em 100)
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A2").Value = "MA"
$wsTx.Range("A3").Value = "XX"
$wsTx.Range("A4").Value = "MT"
